$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# B9: change role from "Requirement Specifier" to "System Analyst "
$ws.Range("B9").Value = "System Analyst "

# Row 12: Gør klar og lav brugertest med HØK / Test Manager
$ws.Range("A12").Value = "Gør klar og lav brugertest med HØK"
$ws.Range("B12").Value = "Test Manager"
$ws.Range("C12").Value = 43887
$ws.Range("D12").Value = 0.43055555555555558
$ws.Range("E12").Value = 0.49305555555555558

# Row 14 F value first, to control shared-string insertion order ("15 min")
$ws.Range("F14").Value = "15 min"

# Row 15: Lav DOM05 beregn markedsføringsbidrag / business-Process Analyst
$ws.Range("A15").Value = "Lav DOM05 beregn markedsføringsbidrag"
$ws.Range("B15").Value = "business-Process Analyst"
$ws.Range("C15").Value = 43887
$ws.Range("D15").Value = 0.57638888888888895
$ws.Range("E15").Value = 0.59027777777777779
$ws.Range("F15").Value = "15 min"

# Row 13: Kundemøde med HØK / Requirement Specifier
$ws.Range("A13").Value = "Kundemøde med HØK"
$ws.Range("B13").Value = "Requirement Specifier"
$ws.Range("C13").Value = 43887
$ws.Range("D13").Value = 0.52083333333333337
$ws.Range("E13").Value = 0.5625
$ws.Range("F13").Value = "45 min"

# Row 14: Lav UC05 beregn markedsføringsbidrag / Requirement Specifier
$ws.Range("A14").Value = "Lav UC05 beregn markedsføringsbidrag"
$ws.Range("B14").Value = "Requirement Specifier"
$ws.Range("C14").Value = 43887
$ws.Range("D14").Value = 0.5625
$ws.Range("E14").Value = 0.57638888888888895

# Update the active cell selection to F16
$ws.Range("F16").Select()
